$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows below had their full contents permuted (rows were reordered by date).
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T are identical for every affected row, so
# only D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) actually need updating to
# reflect the new row order.

$ws.Range("D3").Value = 44533
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 16000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 16500
$ws.Range("S3").Value = 825

$ws.Range("D4").Value = 44792
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("S4").Value = 1075

$ws.Range("D5").Value = 44708
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("S5").Value = 1025

$ws.Range("D6").Value = 44761
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20500
$ws.Range("S6").Value = 1025

$ws.Range("D7").Value = 44320
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 16500
$ws.Range("S7").Value = 825

$ws.Range("D8").Value = 44357
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("S8").Value = 725

$ws.Range("D10").Value = 44893
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 21000
$ws.Range("O10").Value = 22000
$ws.Range("P10").Value = 21625
$ws.Range("S10").Value = 1081
